# Replace the controller board line item: Makerbase MKS Gen_L 2.1 Control Board
# -> Arduino Uno + CNC Shield pack (with its AliExpress link), and reduce the
# Mosfet module IRF520 quantity from 2 to 1 (row 26).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: swap the controller board name + purchase link
$ws.Range("A20").Value = "Arduino Uno + CNC Shield pack"
$ws.Range("D20").Value = "https://www.aliexpress.com/item/1005006431685856.html"

# Row 26: quantity change 2 -> 1 (G26/G28 totals recalc automatically)
$ws.Range("C26").Value = 1

# Move the active selection to D21, matching the saved view state
$ws.Range("D21").Select()
